# Updated cryptos list on Fri Dec 29 20:56:32 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.996.60"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.305.36"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.972"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.68%  "
$ws.Range("D16").Value = "2.652.33"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "2.296.02"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "41.997.76"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.77"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.76"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.90%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -5.84%  "
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.129"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  +11.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  -5.21%  "
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.17"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.22"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.13%  "
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.83"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.21%  "
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "1.565.31"
$ws.Range("E51").Value = "  +0.30%  "
